$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion rates in the daily report text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$text = $wsHoja1.Range("A1").Value2
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 12.99 = 52337.79 pesos"), "✅ 1000 Bs = 13.1 = 54116.35 pesos"
$text = $text -replace [regex]::Escape("✅ 52337.79 pesos = 12.89 = 960.39 Bs"), "✅ 54116.35 pesos = 13.33 = 999.69 Bs"
$wsHoja1.Range("A1").Value2 = $text

# --- Sheet "tasas": update the N10/O10/O12 rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value2 = 76.31999999999999
$wsTasas.Range("O10").Value2 = 4130.16
$wsTasas.Range("O12").Value2 = 75
